$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("D1").Value = "=INDEX(LINEST(C1:C18,B1:B18),1,1)"
